$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"
$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data rows ---
$wsForecast.Cells.Item(2,1).Value = 45109.99999999999
$wsForecast.Cells.Item(2,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(2,2).Value = 593
$wsForecast.Cells.Item(2,3).Value = -62.93950325164867
$wsForecast.Cells.Item(2,4).Value = 1264.344665259001
$wsForecast.Cells.Item(3,1).Value = 45130.99999999999
$wsForecast.Cells.Item(3,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(3,2).Value = 583
$wsForecast.Cells.Item(3,3).Value = -83.02971796729338
$wsForecast.Cells.Item(3,4).Value = 1233.902746761944
$wsForecast.Cells.Item(4,1).Value = 45137.99999999999
$wsForecast.Cells.Item(4,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(4,2).Value = 580
$wsForecast.Cells.Item(4,3).Value = -70.07431571030288
$wsForecast.Cells.Item(4,4).Value = 1197.138249649751
$wsForecast.Cells.Item(5,1).Value = 45151.99999999999
$wsForecast.Cells.Item(5,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(5,2).Value = 573
$wsForecast.Cells.Item(5,3).Value = -108.7435777348074
$wsForecast.Cells.Item(5,4).Value = 1232.346401785927
$wsForecast.Cells.Item(6,1).Value = 45158.99999999999
$wsForecast.Cells.Item(6,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(6,2).Value = 570
$wsForecast.Cells.Item(6,3).Value = -66.28885912713567
$wsForecast.Cells.Item(6,4).Value = 1229.287913963558
$wsForecast.Cells.Item(7,1).Value = 45172.99999999999
$wsForecast.Cells.Item(7,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(7,2).Value = 563
$wsForecast.Cells.Item(7,3).Value = -113.2339786579112
$wsForecast.Cells.Item(7,4).Value = 1231.710553329531
$wsForecast.Cells.Item(8,1).Value = 45186.99999999999
$wsForecast.Cells.Item(8,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(8,2).Value = 557
$wsForecast.Cells.Item(8,3).Value = -112.4263613318796
$wsForecast.Cells.Item(8,4).Value = 1205.303691645433
$wsForecast.Cells.Item(9,1).Value = 45200.99999999999
$wsForecast.Cells.Item(9,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(9,2).Value = 550
$wsForecast.Cells.Item(9,3).Value = -18.50104045965723
$wsForecast.Cells.Item(9,4).Value = 1207.089562546552
$wsForecast.Cells.Item(10,1).Value = 45207.99999999999
$wsForecast.Cells.Item(10,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(10,2).Value = 547
$wsForecast.Cells.Item(10,3).Value = -111.6697959108155
$wsForecast.Cells.Item(10,4).Value = 1199.298565453431
$wsForecast.Cells.Item(11,1).Value = 45214.99999999999
$wsForecast.Cells.Item(11,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(11,2).Value = 544
$wsForecast.Cells.Item(11,3).Value = -118.9476168946281
$wsForecast.Cells.Item(11,4).Value = 1197.284356244044
$wsForecast.Cells.Item(12,1).Value = 45228.99999999999
$wsForecast.Cells.Item(12,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(12,2).Value = 537
$wsForecast.Cells.Item(12,3).Value = -66.11948029605236
$wsForecast.Cells.Item(12,4).Value = 1184.593297876329
$wsForecast.Cells.Item(13,1).Value = 45235.99999999999
$wsForecast.Cells.Item(13,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(13,2).Value = 534
$wsForecast.Cells.Item(13,3).Value = -156.0847081109492
$wsForecast.Cells.Item(13,4).Value = 1204.781940614867
$wsForecast.Cells.Item(14,1).Value = 45256.99999999999
$wsForecast.Cells.Item(14,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(14,2).Value = 524
$wsForecast.Cells.Item(14,3).Value = -116.0181480517459
$wsForecast.Cells.Item(14,4).Value = 1212.49863537655
$wsForecast.Cells.Item(15,1).Value = 45277.99999999999
$wsForecast.Cells.Item(15,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(15,2).Value = 514
$wsForecast.Cells.Item(15,3).Value = -136.9582356833643
$wsForecast.Cells.Item(15,4).Value = 1171.837884718634
$wsForecast.Cells.Item(16,1).Value = 45298.99999999999
$wsForecast.Cells.Item(16,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(16,2).Value = 504
$wsForecast.Cells.Item(16,3).Value = -136.4313365943453
$wsForecast.Cells.Item(16,4).Value = 1146.513248125533
$wsForecast.Cells.Item(17,1).Value = 45312.99999999999
$wsForecast.Cells.Item(17,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(17,2).Value = 498
$wsForecast.Cells.Item(17,3).Value = -171.8071209043906
$wsForecast.Cells.Item(17,4).Value = 1166.33276994776
$wsForecast.Cells.Item(18,1).Value = 45326.99999999999
$wsForecast.Cells.Item(18,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(18,2).Value = 491
$wsForecast.Cells.Item(18,3).Value = -158.6341887479822
$wsForecast.Cells.Item(18,4).Value = 1202.925204447414
$wsForecast.Cells.Item(19,1).Value = 45333.99999999999
$wsForecast.Cells.Item(19,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(19,2).Value = 488
$wsForecast.Cells.Item(19,3).Value = -204.0417059941828
$wsForecast.Cells.Item(19,4).Value = 1191.156221447827
$wsForecast.Cells.Item(20,1).Value = 45361.99999999999
$wsForecast.Cells.Item(20,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(20,2).Value = 475
$wsForecast.Cells.Item(20,3).Value = -207.1105214298225
$wsForecast.Cells.Item(20,4).Value = 1162.930957947634
$wsForecast.Cells.Item(21,1).Value = 45396.99999999999
$wsForecast.Cells.Item(21,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(21,2).Value = 458
$wsForecast.Cells.Item(21,3).Value = -193.5486830990327
$wsForecast.Cells.Item(21,4).Value = 1117.406974419842
$wsForecast.Cells.Item(22,1).Value = 45410.99999999999
$wsForecast.Cells.Item(22,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(22,2).Value = 452
$wsForecast.Cells.Item(22,3).Value = -232.0131004079131
$wsForecast.Cells.Item(22,4).Value = 1168.922455432686
$wsForecast.Cells.Item(23,1).Value = 45424.99999999999
$wsForecast.Cells.Item(23,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(23,2).Value = 445
$wsForecast.Cells.Item(23,3).Value = -169.4695211756835
$wsForecast.Cells.Item(23,4).Value = 1140.837314965748
$wsForecast.Cells.Item(24,1).Value = 45466.99999999999
$wsForecast.Cells.Item(24,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(24,2).Value = 426
$wsForecast.Cells.Item(24,3).Value = -212.060840160649
$wsForecast.Cells.Item(24,4).Value = 1137.088415989238
$wsForecast.Cells.Item(25,1).Value = 45473.99999999999
$wsForecast.Cells.Item(25,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(25,2).Value = 422
$wsForecast.Cells.Item(25,3).Value = -257.4527381016355
$wsForecast.Cells.Item(25,4).Value = 1092.466678709995
$wsForecast.Cells.Item(26,1).Value = 45494.99999999999
$wsForecast.Cells.Item(26,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(26,2).Value = 413
$wsForecast.Cells.Item(26,3).Value = -217.5935804741484
$wsForecast.Cells.Item(26,4).Value = 1087.080664183621
$wsForecast.Cells.Item(27,1).Value = 45501.99999999999
$wsForecast.Cells.Item(27,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(27,2).Value = 409
$wsForecast.Cells.Item(27,3).Value = -248.5586900739114
$wsForecast.Cells.Item(27,4).Value = 1090.214211840718
$wsForecast.Cells.Item(28,1).Value = 45508.99999999999
$wsForecast.Cells.Item(28,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(28,2).Value = 406
$wsForecast.Cells.Item(28,3).Value = -307.8310579176223
$wsForecast.Cells.Item(28,4).Value = 1077.786382384677
$wsForecast.Cells.Item(29,1).Value = 45515.99999999999
$wsForecast.Cells.Item(29,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(29,2).Value = 403
$wsForecast.Cells.Item(29,3).Value = -230.0301926785692
$wsForecast.Cells.Item(29,4).Value = 1050.545829230424
$wsForecast.Cells.Item(30,1).Value = 45522.99999999999
$wsForecast.Cells.Item(30,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(30,2).Value = 399
$wsForecast.Cells.Item(30,3).Value = -283.5692232567208
$wsForecast.Cells.Item(30,4).Value = 1041.098910764783
$wsForecast.Cells.Item(31,1).Value = 45529.99999999999
$wsForecast.Cells.Item(31,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(31,2).Value = 396
$wsForecast.Cells.Item(31,3).Value = -210.5214908765584
$wsForecast.Cells.Item(31,4).Value = 1077.188349278137
$wsForecast.Cells.Item(32,1).Value = 45536.99999999999
$wsForecast.Cells.Item(32,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(32,2).Value = 393
$wsForecast.Cells.Item(32,3).Value = -290.5703355840645
$wsForecast.Cells.Item(32,4).Value = 1076.239405427304
$wsForecast.Cells.Item(33,1).Value = 45543.99999999999
$wsForecast.Cells.Item(33,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(33,2).Value = 390
$wsForecast.Cells.Item(33,3).Value = -286.4939248543017
$wsForecast.Cells.Item(33,4).Value = 1036.453609678717
$wsForecast.Cells.Item(34,1).Value = 45550.99999999999
$wsForecast.Cells.Item(34,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(34,2).Value = 386
$wsForecast.Cells.Item(34,3).Value = -274.7145852753623
$wsForecast.Cells.Item(34,4).Value = 1072.964524835972
$wsForecast.Cells.Item(35,1).Value = 45564.99999999999
$wsForecast.Cells.Item(35,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(35,2).Value = 380
$wsForecast.Cells.Item(35,3).Value = -283.6841040174029
$wsForecast.Cells.Item(35,4).Value = 1047.237364493894
$wsForecast.Cells.Item(36,1).Value = 45571.99999999999
$wsForecast.Cells.Item(36,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(36,2).Value = 377
$wsForecast.Cells.Item(36,3).Value = -299.4421041123945
$wsForecast.Cells.Item(36,4).Value = 1030.031172517364
$wsForecast.Cells.Item(37,1).Value = 45578.99999999999
$wsForecast.Cells.Item(37,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(37,2).Value = 373
$wsForecast.Cells.Item(37,3).Value = -255.3274411018787
$wsForecast.Cells.Item(37,4).Value = 1044.594688545899
$wsForecast.Cells.Item(38,1).Value = 45592.99999999999
$wsForecast.Cells.Item(38,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(38,2).Value = 367
$wsForecast.Cells.Item(38,3).Value = -301.3083450900838
$wsForecast.Cells.Item(38,4).Value = 998.4107204421884
$wsForecast.Cells.Item(39,1).Value = 45599.99999999999
$wsForecast.Cells.Item(39,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(39,2).Value = 363
$wsForecast.Cells.Item(39,3).Value = -275.5098684544466
$wsForecast.Cells.Item(39,4).Value = 1039.141649326198
$wsForecast.Cells.Item(40,1).Value = 45627.99999999999
$wsForecast.Cells.Item(40,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(40,2).Value = 350
$wsForecast.Cells.Item(40,3).Value = -321.5655221517194
$wsForecast.Cells.Item(40,4).Value = 1074.39391302218
$wsForecast.Cells.Item(41,1).Value = 45634.99999999999
$wsForecast.Cells.Item(41,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(41,2).Value = 347
$wsForecast.Cells.Item(41,3).Value = -310.5084744341403
$wsForecast.Cells.Item(41,4).Value = 1008.75894036113
$wsForecast.Cells.Item(42,1).Value = 45641.99999999999
$wsForecast.Cells.Item(42,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(42,2).Value = 344
$wsForecast.Cells.Item(42,3).Value = -344.3780098720625
$wsForecast.Cells.Item(42,4).Value = 982.0988011955405
$wsForecast.Cells.Item(43,1).Value = 45648.99999999999
$wsForecast.Cells.Item(43,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(43,2).Value = 340
$wsForecast.Cells.Item(43,3).Value = -312.2170545491732
$wsForecast.Cells.Item(43,4).Value = 1028.738249804615
$wsForecast.Cells.Item(44,1).Value = 45655.99999999999
$wsForecast.Cells.Item(44,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(44,2).Value = 337
$wsForecast.Cells.Item(44,3).Value = -348.1851873725814
$wsForecast.Cells.Item(44,4).Value = 1021.420759687253
$wsForecast.Cells.Item(45,1).Value = 45662.99999999999
$wsForecast.Cells.Item(45,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(45,2).Value = 334
$wsForecast.Cells.Item(45,3).Value = -325.7023739074555
$wsForecast.Cells.Item(45,4).Value = 1024.416325008792
$wsForecast.Cells.Item(46,1).Value = 45669.99999999999
$wsForecast.Cells.Item(46,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(46,2).Value = 331
$wsForecast.Cells.Item(46,3).Value = -277.8196854055672
$wsForecast.Cells.Item(46,4).Value = 954.1532434876351
$wsForecast.Cells.Item(47,1).Value = 45676.99999999999
$wsForecast.Cells.Item(47,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(47,2).Value = 327
$wsForecast.Cells.Item(47,3).Value = -323.2163249094673
$wsForecast.Cells.Item(47,4).Value = 1003.213446174896
$wsForecast.Cells.Item(48,1).Value = 45683.99999999999
$wsForecast.Cells.Item(48,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Cells.Item(48,2).Value = 324
$wsForecast.Cells.Item(48,3).Value = -304.8445941271431
$wsForecast.Cells.Item(48,4).Value = 979.1538058675351

Write-Host "Edit complete"
